# Apply trade #49 close update to the live trading results workbook.

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.55
$summary.Range("B4").Value = -0.45
$summary.Range("B6").Value = 49
$summary.Range("B8").Value = 21
$summary.Range("B9").Value = 36.73

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.55
$status.Range("D4").Value = 49
$status.Range("E4").Value = -0.45
$status.Range("F4").Value = -0.45
$status.Range("G4").Value = 36.73

# --- New trade row data (trade #49, closed 2026-02-17 08:39:25) ---
$rowData = @{
    A = 49
    B = "2026-02-17"
    C = "08:39:25"
    D = "MarketMaking"
    E = "UP"
    F = 0.86
    G = 0.84
    H = "CLOSED"
    I = -2.3256
    J = -0.02
    K = 99.55
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

$sheetsWithTrades = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetsWithTrades) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A50").Value = $rowData.A

    # Force Date/Time columns to stay as plain text (matches source inlineStr),
    # otherwise Excel auto-parses "2026-02-17" / "08:39:25" as date/time serials.
    $ws.Range("B50").NumberFormat = "@"
    $ws.Range("B50").Value = $rowData.B
    $ws.Range("C50").NumberFormat = "@"
    $ws.Range("C50").Value = $rowData.C

    $ws.Range("D50").Value = $rowData.D
    $ws.Range("E50").Value = $rowData.E
    $ws.Range("F50").Value = $rowData.F
    $ws.Range("G50").Value = $rowData.G
    $ws.Range("H50").Value = $rowData.H
    $ws.Range("I50").Value = $rowData.I
    $ws.Range("J50").Value = $rowData.J
    $ws.Range("K50").Value = $rowData.K
    $ws.Range("L50").Value = $rowData.L
    $ws.Range("M50").Value = $rowData.M
    $ws.Range("N50").Value = $rowData.N
    $ws.Range("O50").Value = $rowData.O
    $ws.Range("P50").Value = $rowData.P
    $ws.Range("Q50").Value = $rowData.Q
}
